$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$notesPage = $s.NotesPage
$notesShape = $notesPage.Shapes.AddPlaceholder(2)

$lines = @(
    'it "calculates a 5% tax before validation" do',
    '    item = Item.new(:name => "Candy", :price => 20)',
    '    item.valid?',
    '    item.tax.should == 1.0',
    '  end',
    '',
    '  it "rejects items who''s price is not numeric" do',
    '    item = Item.new(:name => "Candy", :price => "asdfasdf")',
    '    item.should_not be_valid',
    '  end',
    '',
    'validates_numericality_of :price',
    ''
)

$notesText = [string]::Join("`n", $lines)
$notesShape.TextFrame.TextRange.Text = $notesText
